# Update "paises" (countries) COVID-19 stats sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp header
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 15:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1552140
$ws.Range("C4").Value = 1846
$ws.Range("E4").Value = 1101171
$ws.Range("G4").Value = 82
$ws.Range("H4").Value = 92063

# Row 7 - Brasil
$ws.Range("B7").Value = 257396
$ws.Range("C7").Value = 2028
$ws.Range("E7").Value = 139996
$ws.Range("G7").Value = 88
$ws.Range("H7").Value = 16941

# Row 53 - Noruega
$ws.Range("B53").Value = 8264
$ws.Range("C53").Value = 7
$ws.Range("E53").Value = 7999

# Row 88 - Islandia
$ws.Range("D88").Value = 1789
$ws.Range("E88").Value = 3

# Row 90 - Lituania
$ws.Range("B90").Value = 1562
$ws.Range("C90").Value = 15
$ws.Range("D90").Value = 1025
$ws.Range("E90").Value = 477
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 60

# Row 92 - Republica de Yibuti
$ws.Range("D92").Value = 1033
$ws.Range("E92").Value = 478

# Rows 105/106 - Letonia and Sri Lanka swap order (Sri Lanka now above Letonia)
# with updated Sri Lanka figures, Letonia figures unchanged.
$ws.Range("A105").Value = "Sri Lanka"
$ws.Range("B105").Value = 1020
$ws.Range("C105").Value = 28
$ws.Range("D105").Value = 569
$ws.Range("E105").Value = 442
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 9

$ws.Range("A106").Value = "Letonia"
$ws.Range("B106").Value = 1012
$ws.Range("C106").Value = 3
$ws.Range("D106").Value = 694
$ws.Range("E106").Value = 297
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = 21

# Row 117 - Principado de Andorra
$ws.Range("D117").Value = 628
$ws.Range("E117").Value = 82

# Row 135 - Estado de Palestina
$ws.Range("D135").Value = 346
$ws.Range("E135").Value = 40

# Row 150 - Liberia
$ws.Range("B150").Value = 233
$ws.Range("C150").Value = 4
$ws.Range("D150").Value = 125
$ws.Range("E150").Value = 85
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 23
